$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 597.17645
$ws.Range("I19").Value = 780
$ws.Range("J19").Value = 469.2
$ws.Range("K19").Value = 780
$ws.Range("L19").Value = 469.2
$ws.Range("M19").Value = -605
$ws.Range("N19").Value = -819.2
$ws.Range("H32").Value = 6255311.5
$ws.Range("J32").Value = 9094454
$ws.Range("L32").Value = 9094454
$ws.Range("N32").Value = -9095106
$ws.Range("H62").Value = 4100
$ws.Range("I62").Value = 1200
$ws.Range("K62").Value = 1200
$ws.Range("M62").Value = -576
$ws.Range("H65").Value = 4100
$ws.Range("I65").Value = 1200
$ws.Range("K65").Value = 6000
$ws.Range("M65").Value = -2880
$ws.Range("H98").Value = 585362.8
$ws.Range("I98").Value = 654202.9399999999
$ws.Range("K98").Value = 654202.9399999999
$ws.Range("M98").Value = -652704.9399999999
$ws.Range("H122").Value = 585362.8
$ws.Range("I122").Value = 654202.9399999999
$ws.Range("K122").Value = 1962608.82
$ws.Range("M122").Value = -1960158.82
$ws.Range("H132").Value = 1923.1333
$ws.Range("I132").Value = 1895.4546
$ws.Range("J132").Value = 1999.25
$ws.Range("K132").Value = 5686.3638
$ws.Range("L132").Value = 5997.75
$ws.Range("M132").Value = -3156.3638
$ws.Range("N132").Value = -11057.75
$ws.Range("H135").Value = 986.86664
$ws.Range("I135").Value = 807.3929000000001
$ws.Range("J135").Value = 3499.5
$ws.Range("K135").Value = 7266.5361
$ws.Range("L135").Value = 31495.5
$ws.Range("M135").Value = -4731.5361
$ws.Range("N135").Value = -36565.5
$ws.Range("H137").Value = 3289.3408
$ws.Range("I137").Value = 2003.9474
$ws.Range("J137").Value = 4266.24
$ws.Range("K137").Value = 6011.8422
$ws.Range("L137").Value = 12798.72
$ws.Range("M137").Value = -3461.8422
$ws.Range("N137").Value = -17898.72
$ws.Range("H138").Value = 3023.8572
$ws.Range("J138").Value = 3346.2576
$ws.Range("L138").Value = 10038.7728
$ws.Range("N138").Value = -20318.7728
$ws.Range("H141").Value = 3544.4827
$ws.Range("I141").Value = 3599.6428
$ws.Range("K141").Value = 10798.9284
$ws.Range("M141").Value = -5618.928400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14440.348
$ws.Range("I32").Value = 11268.032
$ws.Range("K32").Value = 11268.032
$ws.Range("M32").Value = -10981.032
$ws.Range("H74").Value = 2190.9429
$ws.Range("J74").Value = 3854.8
$ws.Range("L74").Value = 3854.8
$ws.Range("N74").Value = -5602.8
$ws.Range("H77").Value = 2190.9429
$ws.Range("J77").Value = 3854.8
$ws.Range("L77").Value = 19274
$ws.Range("N77").Value = -28010
$ws.Range("H132").Value = 5371.316
$ws.Range("I132").Value = 4865.5557
$ws.Range("K132").Value = 14596.6671
$ws.Range("M132").Value = -12066.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 13250
$ws.Range("J81").Value = 13250
$ws.Range("L81").Value = 13250
$ws.Range("N81").Value = -15372
$ws.Range("H84").Value = 13250
$ws.Range("J84").Value = 13250
$ws.Range("L84").Value = 39750
$ws.Range("N84").Value = -50358
$ws.Range("H86").Value = 1907.4
$ws.Range("I86").Value = 2032.5385
$ws.Range("J86").Value = 1675
$ws.Range("K86").Value = 2032.5385
$ws.Range("L86").Value = 1675
$ws.Range("M86").Value = -909.5385000000001
$ws.Range("N86").Value = -3921
$ws.Range("H89").Value = 1907.4
$ws.Range("I89").Value = 2032.5385
$ws.Range("J89").Value = 1675
$ws.Range("K89").Value = 10162.6925
$ws.Range("L89").Value = 8375
$ws.Range("M89").Value = -4546.692500000001
$ws.Range("N89").Value = -19607
$ws.Range("H94").Value = 6251882.5
$ws.Range("I94").Value = 1394.1666
$ws.Range("K94").Value = 1394.1666
$ws.Range("M94").Value = -943.1666
$ws.Range("H122").Value = 49999.434
$ws.Range("J122").Value = 49999.434
$ws.Range("L122").Value = 49999.434
$ws.Range("N122").Value = -59799.434

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4010.257
$ws.Range("I31").Value = 3040.6875
$ws.Range("J31").Value = 4826.737
$ws.Range("K31").Value = 3040.6875
$ws.Range("L31").Value = 4826.737
$ws.Range("M31").Value = -2745.6875
$ws.Range("N31").Value = -5416.737
$ws.Range("H34").Value = 4010.257
$ws.Range("I34").Value = 3040.6875
$ws.Range("J34").Value = 4826.737
$ws.Range("K34").Value = 3040.6875
$ws.Range("L34").Value = 4826.737
$ws.Range("M34").Value = -2838.6875
$ws.Range("N34").Value = -5230.737

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3425.6
$ws.Range("I3").Value = 3425.6
$ws.Range("K3").Value = 10276.8
$ws.Range("M3").Value = -10164.8
$ws.Range("H5").Value = 1116.0322
$ws.Range("J5").Value = 2181.9092
$ws.Range("L5").Value = 6545.7276
$ws.Range("N5").Value = -6769.7276
$ws.Range("H56").Value = 7763.625
$ws.Range("I56").Value = 7763.625
$ws.Range("K56").Value = 7763.625
$ws.Range("M56").Value = -7233.625
$ws.Range("H113").Value = 1791.6666
$ws.Range("I113").Value = 1554.25
$ws.Range("J113").Value = 1910.375
$ws.Range("K113").Value = 4662.75
$ws.Range("L113").Value = 5731.125
$ws.Range("M113").Value = -2492.75
$ws.Range("N113").Value = -10071.125
$ws.Range("H133").Value = 2997.5
$ws.Range("I133").Value = 2997.5
$ws.Range("K133").Value = 8992.5
$ws.Range("M133").Value = -3932.5
$ws.Range("H135").Value = 1116.0322
$ws.Range("J135").Value = 2181.9092
$ws.Range("L135").Value = 19637.1828
$ws.Range("N135").Value = -24707.1828
$ws.Range("H136").Value = 3986.7778
$ws.Range("I136").Value = 1900
$ws.Range("K136").Value = 5700
$ws.Range("M136").Value = -600
$ws.Range("H138").Value = 7411.909
$ws.Range("I138").Value = 5994.3335
$ws.Range("K138").Value = 17983.0005
$ws.Range("M138").Value = -12843.0005
$ws.Range("H139").Value = 7170.4595
$ws.Range("I139").Value = 4639
$ws.Range("J139").Value = 8541.666999999999
$ws.Range("K139").Value = 13917
$ws.Range("L139").Value = 25625.001
$ws.Range("M139").Value = -8777
$ws.Range("N139").Value = -35905.001
$ws.Range("H140").Value = 1213.0834
$ws.Range("I140").Value = 756
$ws.Range("J140").Value = 3498.5
$ws.Range("K140").Value = 2268
$ws.Range("L140").Value = 10495.5
$ws.Range("M140").Value = 2912
$ws.Range("N140").Value = -20855.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4288.2666
$ws.Range("I122").Value = 2369.3333
$ws.Range("J122").Value = 7166.6665
$ws.Range("K122").Value = 7107.999899999999
$ws.Range("L122").Value = 21499.9995
$ws.Range("M122").Value = -4657.999899999999
$ws.Range("N122").Value = -26399.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4701.7144
$ws.Range("I7").Value = 2789.3333
$ws.Range("K7").Value = 2789.3333
$ws.Range("M7").Value = -2677.3333
$ws.Range("H82").Value = 1642.8667
$ws.Range("I82").Value = 1005.125
$ws.Range("J82").Value = 2371.7144
$ws.Range("K82").Value = 1005.125
$ws.Range("L82").Value = 2371.7144
$ws.Range("M82").Value = -644.125
$ws.Range("N82").Value = -3093.7144
$ws.Range("H85").Value = 1642.8667
$ws.Range("I85").Value = 1005.125
$ws.Range("J85").Value = 2371.7144
$ws.Range("K85").Value = 1005.125
$ws.Range("L85").Value = 2371.7144
$ws.Range("M85").Value = 242.875
$ws.Range("N85").Value = -4867.7144
$ws.Range("H126").Value = 4701.7144
$ws.Range("I126").Value = 2789.3333
$ws.Range("K126").Value = 8367.999899999999
$ws.Range("M126").Value = -5897.999899999999
$ws.Range("H132").Value = 4261.884
$ws.Range("I132").Value = 3525.0334
$ws.Range("K132").Value = 10575.1002
$ws.Range("M132").Value = -8045.100199999999
$ws.Range("H136").Value = 5173.0713
$ws.Range("I136").Value = 3277.611
$ws.Range("K136").Value = 9832.832999999999
$ws.Range("M136").Value = -7282.832999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4014.6
$ws.Range("I122").Value = 3007.3572
$ws.Range("K122").Value = 9022.071599999999
$ws.Range("M122").Value = -6572.071599999999

Write-Output "Applied 206 cell updates"